$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.793.77"
$ws.Range("E2").Value = "'  -0.15%  "
$ws.Range("D3").Value = "'2.076.89"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'232.86"
$ws.Range("E5").Value = "'  -0.66%  "
$ws.Range("E6").Value = "'  -0.16%  "
$ws.Range("D7").Value = "'58.45"
$ws.Range("E7").Value = "'  -1.24%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("D9").Value = "'0.392"
$ws.Range("E9").Value = "'  +0.30%  "
$ws.Range("E10").Value = "'  -1.11%  "
$ws.Range("E11").Value = "'  +3.16%  "
$ws.Range("E12").Value = "'  +0.97%  "
$ws.Range("D13").Value = "'2.383.82"
$ws.Range("E13").Value = "'  -0.54%  "
$ws.Range("D14").Value = "'21.05"
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("D15").Value = "'0.781"
$ws.Range("E15").Value = "'  +1.28%  "
$ws.Range("D16").Value = "'5.34"
$ws.Range("E16").Value = "'  +0.49%  "
$ws.Range("D17").Value = "'2.069.91"
$ws.Range("E17").Value = "'  -0.93%  "
$ws.Range("D18").Value = "'37.723.51"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("B19").Value = "'Litecoin"
$ws.Range("C19").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'71.57"
$ws.Range("E19").Value = "'  -0.22%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.11"
$ws.Range("E20").Value = "'  -2.55%  "
$ws.Range("D21").Value = "'0.0₃0842"
$ws.Range("E21").Value = "'  +1.39%  "
$ws.Range("D22").Value = "'229.09"
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("E24").Value = "'  -1.18%  "
$ws.Range("E25").Value = "'  +0.13%  "
$ws.Range("D26").Value = "'9.68"
$ws.Range("E26").Value = "'  +6.71%  "
$ws.Range("D27").Value = "'171.64"
$ws.Range("E27").Value = "'  +0.49%  "
$ws.Range("E28").Value = "'  -0.91%  "
$ws.Range("D29").Value = "'19.41"
$ws.Range("E29").Value = "'  -0.88%  "
$ws.Range("E30").Value = "'  -2.26%  "
$ws.Range("E31").Value = "'  +0.89%  "
$ws.Range("D32").Value = "'4.73"
$ws.Range("E32").Value = "'  +0.52%  "
$ws.Range("D33").Value = "'0.0631"
$ws.Range("E33").Value = "'  -0.16%  "
$ws.Range("E34").Value = "'  -0.96%  "
$ws.Range("E35").Value = "'  -2.57%  "
$ws.Range("D37").Value = "'3.39"
$ws.Range("E37").Value = "'  -3.52%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "'  +0.06%  "
$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "'  -1.37%  "
$ws.Range("E40").Value = "'  +7.26%  "
$ws.Range("D41").Value = "'101.70"
$ws.Range("E41").Value = "'  +2.18%  "
$ws.Range("D42").Value = "'0.0975"
$ws.Range("E42").Value = "'  -1.59%  "
$ws.Range("E43").Value = "'  -0.57%  "
$ws.Range("D44").Value = "'17.00"
$ws.Range("E44").Value = "'  +4.57%  "
$ws.Range("D45").Value = "'1.450.70"
$ws.Range("E45").Value = "'  -0.91%  "
$ws.Range("E46").Value = "'  -1.93%  "
$ws.Range("E47").Value = "'  -1.46%  "
$ws.Range("D48").Value = "'4.10"
$ws.Range("E48").Value = "'  -4.83%  "
$ws.Range("D49").Value = "'7.37"
$ws.Range("E49").Value = "'  -1.59%  "
$ws.Range("E50").Value = "'  -1.50%  "
$ws.Range("E51").Value = "'  -0.48%  "
